$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.227295398712158
$ws.Range("B1").Value = 2.176074028015137
$ws.Range("C1").Value = 5.964845657348633
$ws.Range("D1").Value = 1.185987830162048
$ws.Range("E1").Value = 1.282487869262695
